# Fall 2019 -> updated ("additional updates for twenty twenty") class schedule tweak.
# The only real content change is in column A (the "Date, Rm" column): for the
# class sessions that were originally meeting in Rm 270 / Rm CC, the room is
# dropped from the displayed date string (the room assignment no longer
# applies), leaving just the weekday/date. Rows that met in Rm 330 (or already
# had no room / blank room) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value  = "Tues, Aug 27"
$ws.Range("A3").Value  = "Thurs, Aug 29"
$ws.Range("A5").Value  = "Tues, Sep 03"
$ws.Range("A6").Value  = "Thurs, Sep 05"
$ws.Range("A8").Value  = "Thurs, Sep 12"
$ws.Range("A9").Value  = "Tues, Sep 17"
$ws.Range("A10").Value = "Thurs, Sep 19"
$ws.Range("A12").Value = "Thurs, Sep 26"
$ws.Range("A13").Value = "Tues, Oct 01"
$ws.Range("A17").Value = "Tues, Oct 15"
$ws.Range("A21").Value = "Tues, Oct 29"
$ws.Range("A23").Value = "Tues, Nov 05"
$ws.Range("A27").Value = "Tues, Nov 19"
$ws.Range("A29").Value = "Tues, Nov 26"

# Match the author's final view state: scrolled back to the top with B3
# selected (instead of being scrolled down to row 21 with A30 selected).
[void]$ws.Range("B3").Select()
